# "update test data for new template"
#
# 1) Row 1 (headers): columns B..S get re-labelled to match the new
#    template's column layout (A and T..Z are unchanged).
# 2) Rows 2-10 (data): values are rearranged into the new column layout.
#    Column mapping (old -> new):
#      A -> A (unchanged)    C -> B   E -> C   H -> D   I -> E   J -> F
#      L -> G   K -> H   P -> I   M -> J   S -> O   W -> W (unchanged)
#    Columns K,L,M,P,S end up empty; B,D,F,G,O become populated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) fix header row labels (B1..S1) ---
$headers = @{
    2  = "EEZ"
    3  = "FAO area"
    4  = "layer"
    5  = "sector"
    6  = "catch type"
    7  = "year"
    8  = "taxon name"
    9  = "amount"
    10 = "input type"
    11 = "original country fishing "
    12 = "EEZ sub area"
    13 = "subregional area"
    14 = "province state"
    15 = "original sector"
    16 = "original taxon name"
    17 = "original FAO name"
    18 = "adjustment factor"
    19 = "gear type"
}
foreach ($col in $headers.Keys) {
    $ws.Cells.Item(1, $col).Value2 = $headers[$col]
}

# --- 2) rearrange data rows 2-10 ---
for ($r = 2; $r -le 10; $r++) {
    # capture existing values first (use Value2 so strings round-trip via
    # the shared-string table instead of producing COM variant wrapper text)
    $vA = $ws.Cells.Item($r, 1).Value2   # A
    $vC = $ws.Cells.Item($r, 3).Value2   # C
    $vE = $ws.Cells.Item($r, 5).Value2   # E
    $vH = $ws.Cells.Item($r, 8).Value2   # H
    $vI = $ws.Cells.Item($r, 9).Value2   # I
    $vJ = $ws.Cells.Item($r, 10).Value2  # J
    $vK = $ws.Cells.Item($r, 11).Value2  # K
    $vL = $ws.Cells.Item($r, 12).Value2  # L
    $vM = $ws.Cells.Item($r, 13).Value2  # M
    $vP = $ws.Cells.Item($r, 16).Value2  # P
    $vS = $ws.Cells.Item($r, 19).Value2  # S

    # write values into their new homes
    $ws.Cells.Item($r, 2).Value2  = $vC   # B <- C
    $ws.Cells.Item($r, 3).Value2  = $vE   # C <- E
    $ws.Cells.Item($r, 4).Value2  = $vH   # D <- H
    $ws.Cells.Item($r, 5).Value2  = $vI   # E <- I
    $ws.Cells.Item($r, 6).Value2  = $vJ   # F <- J
    $ws.Cells.Item($r, 7).Value2  = $vL   # G <- L
    $ws.Cells.Item($r, 8).Value2  = $vK   # H <- K
    $ws.Cells.Item($r, 9).Value2  = $vP   # I <- P
    $ws.Cells.Item($r, 10).Value2 = $vM   # J <- M
    $ws.Cells.Item($r, 15).Value2 = $vS   # O <- S

    # clear the now-unused source columns
    $ws.Cells.Item($r, 11).ClearContents() # K
    $ws.Cells.Item($r, 12).ClearContents() # L
    $ws.Cells.Item($r, 13).ClearContents() # M
    $ws.Cells.Item($r, 16).ClearContents() # P
    $ws.Cells.Item($r, 19).ClearContents() # S

    # fix up styles so populated/empty cells use the right xf.
    # Reference cells: column A keeps the row's "populated" style, column N
    # keeps the row's "empty" style throughout this edit, so copy formats
    # from them onto the cells whose populated/empty state just flipped.
    $ws.Cells.Item($r, 1).Copy() | Out-Null                 # A = populated style
    $ws.Range($ws.Cells.Item($r, 2), $ws.Cells.Item($r, 2)).PasteSpecial(-4122) | Out-Null  # B
    $ws.Cells.Item($r, 1).Copy() | Out-Null
    $ws.Range($ws.Cells.Item($r, 4), $ws.Cells.Item($r, 4)).PasteSpecial(-4122) | Out-Null  # D
    $ws.Cells.Item($r, 1).Copy() | Out-Null
    $ws.Range($ws.Cells.Item($r, 6), $ws.Cells.Item($r, 6)).PasteSpecial(-4122) | Out-Null  # F
    $ws.Cells.Item($r, 1).Copy() | Out-Null
    $ws.Range($ws.Cells.Item($r, 7), $ws.Cells.Item($r, 7)).PasteSpecial(-4122) | Out-Null  # G
    $ws.Cells.Item($r, 1).Copy() | Out-Null
    $ws.Range($ws.Cells.Item($r, 15), $ws.Cells.Item($r, 15)).PasteSpecial(-4122) | Out-Null # O

    $ws.Cells.Item($r, 14).Copy() | Out-Null                # N = empty style
    $ws.Range($ws.Cells.Item($r, 11), $ws.Cells.Item($r, 11)).PasteSpecial(-4122) | Out-Null # K
    $ws.Cells.Item($r, 14).Copy() | Out-Null
    $ws.Range($ws.Cells.Item($r, 12), $ws.Cells.Item($r, 12)).PasteSpecial(-4122) | Out-Null # L
    $ws.Cells.Item($r, 14).Copy() | Out-Null
    $ws.Range($ws.Cells.Item($r, 13), $ws.Cells.Item($r, 13)).PasteSpecial(-4122) | Out-Null # M
    $ws.Cells.Item($r, 14).Copy() | Out-Null
    $ws.Range($ws.Cells.Item($r, 16), $ws.Cells.Item($r, 16)).PasteSpecial(-4122) | Out-Null # P
    $ws.Cells.Item($r, 14).Copy() | Out-Null
    $ws.Range($ws.Cells.Item($r, 19), $ws.Cells.Item($r, 19)).PasteSpecial(-4122) | Out-Null # S
}

$excel.CutCopyMode = $false
